# Updated TestData for Portugal Market
#
# 1. Duplicate the "Swiss" tab's sibling "Germany" tab (same row/column
#    layout, same generic product names) to seed a new "Portugal" sheet at
#    the end of the workbook, then overwrite the market name + part number.
# 2. Fix up the Germany tab's saved selection.
# 3. Make "Portugal" the active sheet/tab.

$wb = $excel.ActiveWorkbook

$germany = $wb.Worksheets.Item("Germany")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$germany.Copy($null, $lastSheet)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "Portugal"

$ws.Range("B2").Value = "Portugal Market"
$ws.Range("B4").Value = "NGC-3479/T2436"

$ws.Range("B4:B5").Select()

$germanyWs = $wb.Worksheets.Item("Germany")
$germanyWs.Range("A1:D19").Select()

$ws.Activate()
